# Auto-generated script applying scheduled market-data refresh to the
# Exodus_Profits workbook (currentAveragePrice*, LevePrice*, LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 5400
$ws.Range("J26").Value = 5600
$ws.Range("L26").Value = 5600
$ws.Range("N26").Value = -6288
$ws.Range("H80").Value = 2761.7273
$ws.Range("J80").Value = 3625
$ws.Range("L80").Value = 10875
$ws.Range("N80").Value = -12871
$ws.Range("H83").Value = 2761.7273
$ws.Range("J83").Value = 3625
$ws.Range("L83").Value = 32625
$ws.Range("N83").Value = -42609
$ws.Range("H92").Value = 1037.2941
$ws.Range("I92").Value = 731.7143
$ws.Range("K92").Value = 731.7143
$ws.Range("M92").Value = 516.2857
$ws.Range("H98").Value = 650.8889
$ws.Range("I98").Value = 607.1875
$ws.Range("K98").Value = 607.1875
$ws.Range("M98").Value = 890.8125
$ws.Range("H116").Value = 12039.8
$ws.Range("I116").Value = 12599.75
$ws.Range("K116").Value = 12599.75
$ws.Range("M116").Value = -9157.75
$ws.Range("H122").Value = 650.8889
$ws.Range("I122").Value = 607.1875
$ws.Range("K122").Value = 1821.5625
$ws.Range("M122").Value = 628.4375
$ws.Range("H125").Value = 1992.5
$ws.Range("I125").Value = 1985
$ws.Range("K125").Value = 17865
$ws.Range("M125").Value = -15405
$ws.Range("H131").Value = 2093.3333
$ws.Range("I131").Value = 2093.3333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 6279.999899999999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1239.999899999999
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 1201.3125
$ws.Range("I132").Value = 1194.9788
$ws.Range("K132").Value = 3584.936400000001
$ws.Range("M132").Value = -1054.936400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4013.7144
$ws.Range("I74").Value = 1859.4445
$ws.Range("J74").Value = 7891.4
$ws.Range("K74").Value = 1859.4445
$ws.Range("L74").Value = 7891.4
$ws.Range("M74").Value = -985.4445000000001
$ws.Range("N74").Value = -9639.4
$ws.Range("H77").Value = 4013.7144
$ws.Range("I77").Value = 1859.4445
$ws.Range("J77").Value = 7891.4
$ws.Range("K77").Value = 9297.2225
$ws.Range("L77").Value = 39457
$ws.Range("M77").Value = -4929.2225
$ws.Range("N77").Value = -48193
$ws.Range("H110").Value = 2690
$ws.Range("I110").Value = 2475
$ws.Range("J110").Value = 2833.3333
$ws.Range("K110").Value = 2475
$ws.Range("L110").Value = 2833.3333
$ws.Range("M110").Value = -430
$ws.Range("N110").Value = -6923.3333
$ws.Range("H122").Value = 2116.75
$ws.Range("I122").Value = 2057.1052
$ws.Range("K122").Value = 6171.3156
$ws.Range("M122").Value = -3721.3156
$ws.Range("H132").Value = 3729.875
$ws.Range("I132").Value = 3100.6333
$ws.Range("K132").Value = 9301.8999
$ws.Range("M132").Value = -6771.8999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 698.125
$ws.Range("I5").Value = 564.1667
$ws.Range("J5").Value = 1100
$ws.Range("K5").Value = 564.1667
$ws.Range("L5").Value = 1100
$ws.Range("M5").Value = -451.1667
$ws.Range("N5").Value = -1326
$ws.Range("H20").Value = 461287.28
$ws.Range("I20").Value = 803879
$ws.Range("K20").Value = 803879
$ws.Range("M20").Value = -803632
$ws.Range("H80").Value = 1220.4584
$ws.Range("J80").Value = 1316.5264
$ws.Range("L80").Value = 1316.5264
$ws.Range("N80").Value = -3312.5264
$ws.Range("H83").Value = 1220.4584
$ws.Range("J83").Value = 1316.5264
$ws.Range("L83").Value = 6582.632
$ws.Range("N83").Value = -16566.632
$ws.Range("H94").Value = 1101.871
$ws.Range("I94").Value = 963
$ws.Range("J94").Value = 2398
$ws.Range("K94").Value = 963
$ws.Range("L94").Value = 2398
$ws.Range("M94").Value = -512
$ws.Range("N94").Value = -3300
$ws.Range("H105").Value = 251974.5
$ws.Range("I105").Value = 334900
$ws.Range("J105").Value = 3198
$ws.Range("K105").Value = 334900
$ws.Range("L105").Value = 3198
$ws.Range("M105").Value = -333153
$ws.Range("N105").Value = -6692

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7902658
$ws.Range("I4").Value = 200130
$ws.Range("J4").Value = 15605186
$ws.Range("K4").Value = 200130
$ws.Range("L4").Value = 15605186
$ws.Range("M4").Value = -200018
$ws.Range("N4").Value = -15605410
$ws.Range("H31").Value = 2119.6072
$ws.Range("I31").Value = 1523
$ws.Range("J31").Value = 3379.111
$ws.Range("K31").Value = 1523
$ws.Range("L31").Value = 3379.111
$ws.Range("M31").Value = -1228
$ws.Range("N31").Value = -3969.111
$ws.Range("H34").Value = 2119.6072
$ws.Range("I34").Value = 1523
$ws.Range("J34").Value = 3379.111
$ws.Range("K34").Value = 1523
$ws.Range("L34").Value = 3379.111
$ws.Range("M34").Value = -1321
$ws.Range("N34").Value = -3783.111
$ws.Range("H58").Value = 3015.0881
$ws.Range("I58").Value = 2807.0435
$ws.Range("J58").Value = 3450.0908
$ws.Range("K58").Value = 2807.0435
$ws.Range("L58").Value = 3450.0908
$ws.Range("M58").Value = -2604.0435
$ws.Range("N58").Value = -3856.0908
$ws.Range("H99").Value = 55558304
$ws.Range("J99").Value = 5500
$ws.Range("L99").Value = 5500
$ws.Range("N99").Value = -8496
$ws.Range("H126").Value = 55558304
$ws.Range("J126").Value = 5500
$ws.Range("L126").Value = 16500
$ws.Range("N126").Value = -21440
$ws.Range("H132").Value = 1106459.1
$ws.Range("I132").Value = 1037096.6
$ws.Range("K132").Value = 3111289.8
$ws.Range("M132").Value = -3108759.8
$ws.Range("H134").Value = 2166890.8
$ws.Range("I134").Value = 2749394
$ws.Range("K134").Value = 8248182
$ws.Range("M134").Value = -8245647
$ws.Range("H136").Value = 3015.0881
$ws.Range("I136").Value = 2807.0435
$ws.Range("J136").Value = 3450.0908
$ws.Range("K136").Value = 8421.130500000001
$ws.Range("L136").Value = 10350.2724
$ws.Range("M136").Value = -5871.130500000001
$ws.Range("N136").Value = -15450.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3333711
$ws.Range("I4").Value = 3448662
$ws.Range("K4").Value = 10345986
$ws.Range("M4").Value = -10345874
$ws.Range("H38").Value = 243
$ws.Range("I38").Value = 5.7777777
$ws.Range("J38").Value = 1310.5
$ws.Range("K38").Value = 17.3333331
$ws.Range("L38").Value = 3931.5
$ws.Range("M38").Value = 329.6666669
$ws.Range("N38").Value = -4625.5
$ws.Range("H129").Value = 487
$ws.Range("I129").Value = 605.75
$ws.Range("J129").Value = 249.5
$ws.Range("K129").Value = 1817.25
$ws.Range("L129").Value = 748.5
$ws.Range("M129").Value = 3182.75
$ws.Range("N129").Value = -10748.5
$ws.Range("H133").Value = 6190.625
$ws.Range("I133").Value = 4920.8335
$ws.Range("K133").Value = 14762.5005
$ws.Range("M133").Value = -9702.500499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 197
$ws.Range("I5").Value = 197
$ws.Range("K5").Value = 197
$ws.Range("M5").Value = -85
$ws.Range("H31").Value = 1732.75
$ws.Range("I31").Value = 1730.2858
$ws.Range("J31").Value = 1750
$ws.Range("K31").Value = 1730.2858
$ws.Range("L31").Value = 1750
$ws.Range("M31").Value = -1438.2858
$ws.Range("N31").Value = -2334
$ws.Range("H37").Value = 1732.75
$ws.Range("I37").Value = 1730.2858
$ws.Range("J37").Value = 1750
$ws.Range("K37").Value = 1730.2858
$ws.Range("L37").Value = 1750
$ws.Range("M37").Value = -1453.2858
$ws.Range("N37").Value = -2304
$ws.Range("H70").Value = 6203.2
$ws.Range("I70").Value = 6335.6665
$ws.Range("K70").Value = 6335.6665
$ws.Range("M70").Value = -6065.6665
$ws.Range("H73").Value = 6203.2
$ws.Range("I73").Value = 6335.6665
$ws.Range("K73").Value = 6335.6665
$ws.Range("M73").Value = -5399.6665
$ws.Range("H102").Value = 3333
$ws.Range("I102").Value = 3333
$ws.Range("K102").Value = 3333
$ws.Range("M102").Value = -1711
$ws.Range("H107").Value = 780.05
$ws.Range("J107").Value = 523.1667
$ws.Range("L107").Value = 523.1667
$ws.Range("N107").Value = -4363.1667
$ws.Range("H122").Value = 95796.21000000001
$ws.Range("I122").Value = 125837.89
$ws.Range("J122").Value = 5671.1665
$ws.Range("K122").Value = 377513.67
$ws.Range("L122").Value = 17013.4995
$ws.Range("M122").Value = -375063.67
$ws.Range("N122").Value = -21913.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8105.32
$ws.Range("I7").Value = 12749.091
$ws.Range("K7").Value = 12749.091
$ws.Range("M7").Value = -12637.091
$ws.Range("H16").Value = 2541
$ws.Range("J16").Value = 2389.1667
$ws.Range("L16").Value = 2389.1667
$ws.Range("N16").Value = -2729.1667
$ws.Range("H40").Value = 4717600.5
$ws.Range("I40").Value = 104985.1
$ws.Range("K40").Value = 104985.1
$ws.Range("M40").Value = -104849.1
$ws.Range("H126").Value = 8105.32
$ws.Range("I126").Value = 12749.091
$ws.Range("K126").Value = 38247.273
$ws.Range("M126").Value = -35777.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1484.2
$ws.Range("I107").Value = 1426.8889
$ws.Range("K107").Value = 4280.6667
$ws.Range("M107").Value = -2360.6667
$ws.Range("H113").Value = 5529.4
$ws.Range("I113").Value = 5062.125
$ws.Range("J113").Value = 7398.5
$ws.Range("K113").Value = 15186.375
$ws.Range("L113").Value = 22195.5
$ws.Range("M113").Value = -13016.375
$ws.Range("N113").Value = -26535.5
$ws.Range("H126").Value = 36707.062
$ws.Range("I126").Value = 48955.543
$ws.Range("K126").Value = 146866.629
$ws.Range("M126").Value = -144396.629
